# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# ---- Overview sheet: update per-language status cells ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $status
$overview.Range("F2").Value = $status
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status
$overview.Range("E4").Value = $status
$overview.Range("F4").Value = $status

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $status
$zhcn.Range("C3").Value = $status
$zhcn.Range("C4").Value = $status

$zhcn.Range("I2").Value = "7e68274c-742c-4c7e-843d-37352c360c5d.md"
$zhcn.Range("J2").Value = "7e68274c-742c-4c7e-843d-37352c360c5d.61741a3c004d6de7075fb8a824a9652e72e967fa.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-05 12:38:33"

$zhcn.Range("I3").Value = "da8c2352-9061-4275-9c03-f12791665aa8.yml"
$zhcn.Range("J3").Value = "da8c2352-9061-4275-9c03-f12791665aa8.8a61b96ab902889e85889dee8a34d717a76a6f6b.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-05 12:38:33"

$zhcn.Range("I4").Value = "fe3cfacb-9e31-4cc9-a7e3-e5ab07d67260.yml"
$zhcn.Range("J4").Value = "fe3cfacb-9e31-4cc9-a7e3-e5ab07d67260.4fd7c3f6b040d94fe75b91c321f37e45b43f41b2.zh-cn.xlf"
$zhcn.Range("K4").Value = "2016-09-05 12:38:33"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb49a404e50f80fd500f0f67a8d2625e589ea01a/e2e/7e68274c-742c-4c7e-843d-37352c360c5d.md", "", "", "7e68274c-742c-4c7e-843d-37352c360c5d.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb49a404e50f80fd500f0f67a8d2625e589ea01a/e2e/da8c2352-9061-4275-9c03-f12791665aa8.yml", "", "", "da8c2352-9061-4275-9c03-f12791665aa8.yml")
$zhcn.Hyperlinks.Add($zhcn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb49a404e50f80fd500f0f67a8d2625e589ea01a/e2e/fe3cfacb-9e31-4cc9-a7e3-e5ab07d67260.yml", "", "", "fe3cfacb-9e31-4cc9-a7e3-e5ab07d67260.yml")

$zhcn.Range("I2:I4").Style = "HyperLink"

$zhcn.Columns.Item(3).ColumnWidth = 29.9777050018311
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $status
$dede.Range("C3").Value = $status
$dede.Range("C4").Value = $status

$dede.Range("I2").Value = "7e68274c-742c-4c7e-843d-37352c360c5d.md"
$dede.Range("J2").Value = "7e68274c-742c-4c7e-843d-37352c360c5d.61741a3c004d6de7075fb8a824a9652e72e967fa.de-de.xlf"
$dede.Range("K2").Value = "2016-09-05 12:38:51"

$dede.Range("I3").Value = "da8c2352-9061-4275-9c03-f12791665aa8.yml"
$dede.Range("J3").Value = "da8c2352-9061-4275-9c03-f12791665aa8.8a61b96ab902889e85889dee8a34d717a76a6f6b.de-de.xlf"
$dede.Range("K3").Value = "2016-09-05 12:38:51"

$dede.Range("I4").Value = "fe3cfacb-9e31-4cc9-a7e3-e5ab07d67260.yml"
$dede.Range("J4").Value = "fe3cfacb-9e31-4cc9-a7e3-e5ab07d67260.4fd7c3f6b040d94fe75b91c321f37e45b43f41b2.de-de.xlf"
$dede.Range("K4").Value = "2016-09-05 12:38:51"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb49a404e50f80fd500f0f67a8d2625e589ea01a/e2e/7e68274c-742c-4c7e-843d-37352c360c5d.md", "", "", "7e68274c-742c-4c7e-843d-37352c360c5d.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb49a404e50f80fd500f0f67a8d2625e589ea01a/e2e/da8c2352-9061-4275-9c03-f12791665aa8.yml", "", "", "da8c2352-9061-4275-9c03-f12791665aa8.yml")
$dede.Hyperlinks.Add($dede.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb49a404e50f80fd500f0f67a8d2625e589ea01a/e2e/fe3cfacb-9e31-4cc9-a7e3-e5ab07d67260.yml", "", "", "fe3cfacb-9e31-4cc9-a7e3-e5ab07d67260.yml")

$dede.Range("I2:I4").Style = "HyperLink"

$dede.Columns.Item(3).ColumnWidth = 29.9777050018311
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40

# ---- Overview column widths (status text widened) ----
$overview.Columns.Item(5).ColumnWidth = 29.9777050018311
$overview.Columns.Item(6).ColumnWidth = 29.9777050018311
